# Fruta / hortaliza, semanal
# Insert a new daily-price record as row 279 (pushing the existing rows
# 279-356 down to 280-357), matching the weekly refresh of the
# "Feria Lagunitas de Puerto Montt - Piña" consolidated sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 279:356 down to 280:357, leaving a blank row 279 to fill in.
$ws.Rows.Item(279).Insert()

# Populate the new row 279 with the latest observation.
$ws.Cells.Item(279, 1).Value = 4
$ws.Cells.Item(279, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(279, 3).Value = "Los Lagos"
$ws.Cells.Item(279, 4).Value = 44932
$ws.Cells.Item(279, 5).Value = 10
$ws.Cells.Item(279, 6).Value = "Fruta"
$ws.Cells.Item(279, 7).Value = 100108
$ws.Cells.Item(279, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(279, 9).Value = 100108005
$ws.Cells.Item(279, 10).Value = "Piña"
$ws.Cells.Item(279, 11).Value = "Caramelo"
$ws.Cells.Item(279, 12).Value = "Segunda"
$ws.Cells.Item(279, 13).Value = 160
$ws.Cells.Item(279, 14).Value = 22000
$ws.Cells.Item(279, 15).Value = 22000
$ws.Cells.Item(279, 16).Value = 22000
$ws.Cells.Item(279, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(279, 18).Value = "Ecuador"
$ws.Cells.Item(279, 19).Value = 1571
$ws.Cells.Item(279, 20).Value = 14
